$wb = $excel.ActiveWorkbook

# --- Water sheet: update labels/values for transport & specific cost (euros -> USD) ---
$wsWater = $wb.Worksheets.Item("Water")
$wsWater.Range("A4").Value = "Water transport cost (USD/100 km/m3)"
$wsWater.Range("B4").Value = 0.06
$wsWater.Range("A5").Value = "Water specific cost (USD/m3)"
$wsWater.Range("B5").Value = 0.247

# --- Infra sheet: replace computed formulas with plain, updated numbers ---
$wsInfra = $wb.Worksheets.Item("Infra")
$wsInfra.Range("B3").Value = 515563
$wsInfra.Range("C3").Value = 5155
$wsInfra.Range("B4").Value = 515563
$wsInfra.Range("C4").Value = 5155

# --- Update active sheet / selections to match the latest saved state ---
[void]$wsWater.Activate()
[void]$wsWater.Range("B5").Select()

[void]$wsInfra.Activate()
[void]$wsInfra.Range("G4").Select()
